$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted at row 96; every existing record
# from row 96 down to the last row (242) shifts down by one row (so the
# former row 242 becomes row 243).
$ws.Rows(96).Insert()

# Populate the freshly inserted row 96 with the new weekly record. The
# descriptive columns (Mercado/Region/Categoria/etc.) are identical to
# every other row in this sheet; only the date (D) and volume (J) are new
# values, with the price columns (K/L/M/P) at the series' usual defaults.
$row = 96
$ws.Cells.Item($row, 1).Value2 = 3                          # A - Mercado ID
$ws.Cells.Item($row, 2).Value2 = "Femacal de La Calera"      # B - Mercado
$ws.Cells.Item($row, 3).Value2 = "Coquimbo"                  # C - Region
$ws.Cells.Item($row, 4).Value2 = 44571                       # D - Fecha
$ws.Cells.Item($row, 5).Value2 = 5                           # E - Codreg
$ws.Cells.Item($row, 6).Value2 = 100112039                   # F - Categoria ID
$ws.Cells.Item($row, 7).Value2 = "Ciboulette"                 # G - Categoria
$ws.Cells.Item($row, 8).Value2 = "Sin especificar"            # H - Variedad
$ws.Cells.Item($row, 9).Value2 = "Primera"                    # I - Calidad
$ws.Cells.Item($row, 10).Value2 = 190                         # J - Volumen
$ws.Cells.Item($row, 11).Value2 = 1500                        # K - Precio minimo
$ws.Cells.Item($row, 12).Value2 = 1500                        # L - Precio maximo
$ws.Cells.Item($row, 13).Value2 = 1500                        # M - Precio promedio ponderado
$ws.Cells.Item($row, 14).Value2 = "$/docena de atados"        # N - Unidad de comercializacion
$ws.Cells.Item($row, 15).Value2 = "Provincia de Quillota"      # O - Origen
$ws.Cells.Item($row, 16).Value2 = 500                          # P - Precio $/Kg
$ws.Cells.Item($row, 17).Value2 = 3                            # Q - Kg o Unidades
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"                  # R - Clasificacion

Write-Host "row inserted"
